# The commit swaps the presentation's theme: the deck's live theme (the
# theme part used by the slide master, exposed here through
# $p.SlideMaster.Theme) was "Integral" and becomes the plain "Office
# Theme" colour palette (the content that used to live in the unused
# theme1.xml part, which only the notes master referenced and which this
# COM object model does not expose separately).
#
# The font scheme (Arial throughout) and the format scheme (fills/lines/
# effects) are byte-identical between the two themes, so the only real
# difference to reproduce through the object model is the 12-slot colour
# scheme: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# Office Theme colour scheme (RGB packed as &HBBGGRR, PowerPoint's native
# long-colour encoding used by the ColorFormat.RGB property).
$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
